$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert new row at position 4 (pushes old rows 4..18 down to 5..19) ---
$ws.Rows.Item(4).Insert()
$ws.Cells.Item(4, 1).Value = 11
$ws.Cells.Item(4, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(4, 3).Value = "Bíobío"
$ws.Cells.Item(4, 4).Value = 44428
$ws.Cells.Item(4, 5).Value = 8
$ws.Cells.Item(4, 6).Value = 100112013
$ws.Cells.Item(4, 7).Value = "Alcachofa"
$ws.Cells.Item(4, 8).Value = "Española"
$ws.Cells.Item(4, 9).Value = "Primera"
$ws.Cells.Item(4, 10).Value = 100
$ws.Cells.Item(4, 11).Value = 14000
$ws.Cells.Item(4, 12).Value = 15000
$ws.Cells.Item(4, 13).Value = 14500
$ws.Cells.Item(4, 14).Value = "$/caja 30 unidades"
$ws.Cells.Item(4, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(4, 16).Value = 483
$ws.Cells.Item(4, 17).Value = 30
$ws.Cells.Item(4, 18).Value = "Hortaliza"

# --- Insert new row at position 13 (pushes rows 13..19 down to 14..20) ---
$ws.Rows.Item(13).Insert()
$ws.Cells.Item(13, 1).Value = 11
$ws.Cells.Item(13, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(13, 3).Value = "Bíobío"
$ws.Cells.Item(13, 4).Value = 44435
$ws.Cells.Item(13, 5).Value = 8
$ws.Cells.Item(13, 6).Value = 100112013
$ws.Cells.Item(13, 7).Value = "Alcachofa"
$ws.Cells.Item(13, 8).Value = "Argentina(o)"
$ws.Cells.Item(13, 9).Value = "Primera"
$ws.Cells.Item(13, 10).Value = 100
$ws.Cells.Item(13, 11).Value = 14000
$ws.Cells.Item(13, 12).Value = 15000
$ws.Cells.Item(13, 13).Value = 14500
$ws.Cells.Item(13, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(13, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(13, 16).Value = 290
$ws.Cells.Item(13, 17).Value = 50
$ws.Cells.Item(13, 18).Value = "Hortaliza"

# --- Insert new row at position 20 (pushes old last row 19 down to 21) ---
$ws.Rows.Item(20).Insert()
$ws.Cells.Item(20, 1).Value = 11
$ws.Cells.Item(20, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(20, 3).Value = "Bíobío"
$ws.Cells.Item(20, 4).Value = 44433
$ws.Cells.Item(20, 5).Value = 8
$ws.Cells.Item(20, 6).Value = 100112013
$ws.Cells.Item(20, 7).Value = "Alcachofa"
$ws.Cells.Item(20, 8).Value = "Argentina(o)"
$ws.Cells.Item(20, 9).Value = "Primera"
$ws.Cells.Item(20, 10).Value = 100
$ws.Cells.Item(20, 11).Value = 14000
$ws.Cells.Item(20, 12).Value = 15000
$ws.Cells.Item(20, 13).Value = 14500
$ws.Cells.Item(20, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(20, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(20, 16).Value = 290
$ws.Cells.Item(20, 17).Value = 50
$ws.Cells.Item(20, 18).Value = "Hortaliza"
